$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 908970.3
$ws.Range("J33").Value = 4531.6665
$ws.Range("L33").Value = 4531.6665
$ws.Range("N33").Value = -4989.6665

$ws.Range("H64").Value = 6206.125
$ws.Range("I64").Value = 3673.75
$ws.Range("K64").Value = 3673.75
$ws.Range("M64").Value = -3425.75

$ws.Range("H67").Value = 6206.125
$ws.Range("I67").Value = 3673.75
$ws.Range("K67").Value = 3673.75
$ws.Range("M67").Value = -2815.75

$ws.Range("H70").Value = 36800.8
$ws.Range("I70").Value = 55668
$ws.Range("J70").Value = 8500
$ws.Range("K70").Value = 167004
$ws.Range("L70").Value = 25500
$ws.Range("M70").Value = -166734
$ws.Range("N70").Value = -26040

$ws.Range("H73").Value = 36800.8
$ws.Range("I73").Value = 55668
$ws.Range("J73").Value = 8500
$ws.Range("K73").Value = 167004
$ws.Range("L73").Value = 25500
$ws.Range("M73").Value = -166068
$ws.Range("N73").Value = -27372

$ws.Range("H121").Value = 4471.143
$ws.Range("J121").Value = 4471.143
$ws.Range("L121").Value = 13413.429
$ws.Range("N121").Value = -16907.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 196539.06
$ws.Range("I32").Value = 195434.9
$ws.Range("K32").Value = 195434.9
$ws.Range("M32").Value = -195147.9

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H110").Value = 2412.6667
$ws.Range("I110").Value = 2631.5
$ws.Range("J110").Value = 1975
$ws.Range("K110").Value = 2631.5
$ws.Range("L110").Value = 1975
$ws.Range("M110").Value = -586.5
$ws.Range("N110").Value = -6065

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 900000000
$ws.Range("J19").Value = 900000000
$ws.Range("L19").Value = 900000000
$ws.Range("N19").Value = -900000346

$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

$ws.Range("H82").Value = 18722.334
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 18722.334
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H86").Value = 3375.4
$ws.Range("J86").Value = 3922.5
$ws.Range("L86").Value = 3922.5
$ws.Range("N86").Value = -6168.5

$ws.Range("H89").Value = 3375.4
$ws.Range("J89").Value = 3922.5
$ws.Range("L89").Value = 19612.5
$ws.Range("N89").Value = -30844.5

$ws.Range("H105").Value = 5265011
$ws.Range("I105").Value = 7694131.5
$ws.Range("K105").Value = 7694131.5
$ws.Range("M105").Value = -7692384.5

$ws.Range("H107").Value = 11250
$ws.Range("I107").Value = 20000
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 20000
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -18080
$ws.Range("N107").Value = -6340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2183.2727
$ws.Range("I31").Value = 1957.2
$ws.Range("K31").Value = 1957.2
$ws.Range("M31").Value = -1662.2

$ws.Range("H34").Value = 2183.2727
$ws.Range("I34").Value = 1957.2
$ws.Range("K34").Value = 1957.2
$ws.Range("M34").Value = -1755.2

$ws.Range("H86").Value = 7206.0625
$ws.Range("I86").Value = 7432.25
$ws.Range("J86").Value = 6979.875
$ws.Range("K86").Value = 7432.25
$ws.Range("L86").Value = 6979.875
$ws.Range("M86").Value = -6309.25
$ws.Range("N86").Value = -9225.875

$ws.Range("H89").Value = 7206.0625
$ws.Range("I89").Value = 7432.25
$ws.Range("J89").Value = 6979.875
$ws.Range("K89").Value = 37161.25
$ws.Range("L89").Value = 34899.375
$ws.Range("M89").Value = -31545.25
$ws.Range("N89").Value = -46131.375

$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2065.6667
$ws.Range("J80").Value = 1700
$ws.Range("L80").Value = 5100
$ws.Range("N80").Value = -6972

$ws.Range("H83").Value = 2065.6667
$ws.Range("J83").Value = 1700
$ws.Range("L83").Value = 15300
$ws.Range("N83").Value = -24660

$ws.Range("H107").Value = 1392.3077
$ws.Range("I107").Value = 1500
$ws.Range("K107").Value = 4500
$ws.Range("M107").Value = -2580

$ws.Range("H131").Value = 85142.336
$ws.Range("I131").Value = 1541.8182
$ws.Range("J131").Value = 155881.23
$ws.Range("K131").Value = 4625.4546
$ws.Range("L131").Value = 467643.6900000001
$ws.Range("M131").Value = 414.5454
$ws.Range("N131").Value = -477723.6900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12170.77
$ws.Range("I80").Value = 15642
$ws.Range("J80").Value = 6616.8
$ws.Range("K80").Value = 15642
$ws.Range("L80").Value = 6616.8
$ws.Range("M80").Value = -14644
$ws.Range("N80").Value = -8612.799999999999

$ws.Range("H83").Value = 12170.77
$ws.Range("I83").Value = 15642
$ws.Range("J83").Value = 6616.8
$ws.Range("K83").Value = 78210
$ws.Range("L83").Value = 33084
$ws.Range("M83").Value = -73218
$ws.Range("N83").Value = -43068

$ws.Range("H93").Value = 43417
$ws.Range("J93").Value = 43417
$ws.Range("L93").Value = 43417
$ws.Range("N93").Value = -47161

$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 4000
$ws.Range("K102").Value = 4000
$ws.Range("M102").Value = -2378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28524.916
$ws.Range("I7").Value = 36088.777
$ws.Range("J7").Value = 5833.3335
$ws.Range("K7").Value = 36088.777
$ws.Range("L7").Value = 5833.3335
$ws.Range("M7").Value = -35976.777
$ws.Range("N7").Value = -6057.3335

$ws.Range("H40").Value = 3500.6
$ws.Range("I40").Value = 2429.4285
$ws.Range("K40").Value = 2429.4285
$ws.Range("M40").Value = -2293.4285

$ws.Range("H46").Value = 44811.6
$ws.Range("I46").Value = 72186
$ws.Range("J46").Value = 3750
$ws.Range("K46").Value = 72186
$ws.Range("L46").Value = 3750
$ws.Range("M46").Value = -71998
$ws.Range("N46").Value = -4126

$ws.Range("H126").Value = 28524.916
$ws.Range("I126").Value = 36088.777
$ws.Range("J126").Value = 5833.3335
$ws.Range("K126").Value = 108266.331
$ws.Range("L126").Value = 17500.0005
$ws.Range("M126").Value = -105796.331
$ws.Range("N126").Value = -22440.0005

$ws.Range("H132").Value = 13408.7
$ws.Range("I132").Value = 26997
$ws.Range("K132").Value = 80991
$ws.Range("M132").Value = -78461

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40615
$ws.Range("J105").Value = 40615
$ws.Range("L105").Value = 40615
$ws.Range("N105").Value = -47603

$ws.Range("H122").Value = 922.4375
$ws.Range("I122").Value = 922.4375
$ws.Range("K122").Value = 2767.3125
$ws.Range("M122").Value = -317.3125
